$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 ("Exhibition")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10615
$ws1.Range("G3").Value = "不可售"
$ws1.Range("G4").Value = "不可售"
$ws1.Range("G5").Value = "不可售"
$ws1.Range("F6").Value = 807
$ws1.Range("F9").Value = 1151
$ws1.Range("F10").Value = 3337
$ws1.Range("F11").Value = 2465
$ws1.Range("F13").Value = 2272
$ws1.Range("F15").Value = 1943
$ws1.Range("F16").Value = 486
$ws1.Range("F18").Value = 615
$ws1.Range("F20").Value = 276
$ws1.Range("F21").Value = 16
$ws1.Range("F22").Value = 38
$ws1.Range("F23").Value = 255
$ws1.Range("F24").Value = 57
$ws1.Range("F25").Value = 407
$ws1.Range("F26").Value = 18
$ws1.Range("F28").Value = 435
$ws1.Range("F29").Value = 639
$ws1.Range("F30").Value = 65
$ws1.Range("F32").Value = 332
$ws1.Range("F33").Value = 27
$ws1.Range("F34").Value = 1595
$ws1.Range("F35").Value = 684
$ws1.Range("F36").Value = 684
$ws1.Range("F37").Value = 1846
$ws1.Range("F38").Value = 187
$ws1.Range("F39").Value = 481
$ws1.Range("F40").Value = 69
$ws1.Range("F41").Value = 518
$ws1.Range("F42").Value = 1141
$ws1.Range("F44").Value = 389

# Sheet 2: 演出 ("Performance")
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 10

# Sheet 4: 全部类型 ("All Types")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10615
$ws4.Range("G3").Value = "不可售"
$ws4.Range("G4").Value = "不可售"
$ws4.Range("G5").Value = "不可售"
$ws4.Range("F6").Value = 807
$ws4.Range("F8").Value = 1151
$ws4.Range("F9").Value = 3337
$ws4.Range("F10").Value = 2465
$ws4.Range("F11").Value = 2272
$ws4.Range("F12").Value = 1943
$ws4.Range("F13").Value = 486
$ws4.Range("F15").Value = 615
$ws4.Range("F17").Value = 276
$ws4.Range("F18").Value = 16
$ws4.Range("F19").Value = 38
$ws4.Range("F20").Value = 255
$ws4.Range("F21").Value = 57
$ws4.Range("F22").Value = 407
$ws4.Range("F23").Value = 18
$ws4.Range("F25").Value = 435
$ws4.Range("F26").Value = 639
$ws4.Range("F27").Value = 65
$ws4.Range("F32").Value = 332
$ws4.Range("F33").Value = 27
$ws4.Range("F34").Value = 1595
$ws4.Range("F35").Value = 684
$ws4.Range("F37").Value = 684
$ws4.Range("F38").Value = 1846
$ws4.Range("F39").Value = 187
$ws4.Range("F43").Value = 481
$ws4.Range("F44").Value = 69
$ws4.Range("F45").Value = 518
$ws4.Range("F46").Value = 1141
$ws4.Range("F48").Value = 389
$ws4.Range("F49").Value = 10
